# Disable "smart quotes" autocorrect so straight quotes in our new text survive.
$word.Options.AutoFormatAsYouTypeReplaceQuotes = $false
$word.Options.AutoFormatReplaceQuotes = $false

$d = $word.ActiveDocument

# --- 1. Insert a new "Meta description" paragraph right after the Heading1 title. ---
$titlePara = $d.Paragraphs.Item(1)
$titlePara.Range.InsertParagraphAfter() | Out-Null

$metaPara = $d.Paragraphs.Item(2)
$metaPara.Style = "Normal"

$metaText = "Meta description: Read our unbiased review of Always Hot Deluxe, a classic slot machine with a non-progressive jackpot of up to 60,000 coins. Play for free and learn more!"
$metaPara.Range.Text = $metaText

$metaStart = $metaPara.Range.Start
$labelLen = "Meta description".Length
$boldRange = $d.Range($metaStart, $metaStart + $labelLen)
$boldRange.Font.Bold = 1

# --- 2. Remove the duplicate bold title paragraph near the end of the document, and
#        replace the following italic paragraph's text with the image-prompt copy. ---
$oldTitleText = "Play Always Hot Deluxe Slot Game for Free - Review"
$oldMetaText = "Read our unbiased review of Always Hot Deluxe, a classic slot machine with a non-progressive jackpot of up to 60,000 coins. Play for free and learn more!"
$imagePromptText = 'Create a feature image that captures the essence of Always Hot Deluxe. The image should be in cartoon style and feature a happy Maya warrior with glasses. The Maya warrior should be holding a slot machine with flames and hot red fruits bursting out of it. The background should be filled with flames and the words "Always Hot Deluxe" should be written in bold and fiery letters.'

# Walk backwards (skipping the very first, Heading-1-styled paragraph) so the
# duplicate plain-body title/description pair near the end of the document are
# the only ones matched.
for ($i = $d.Paragraphs.Count; $i -ge 2; $i--) {
    $para = $d.Paragraphs.Item($i)
    $paraText = $para.Range.Text
    $trimmed = $paraText.TrimEnd("`r", "`n", "`x07")

    if ($trimmed -eq $oldTitleText) {
        $para.Range.Delete() | Out-Null
    }
    elseif ($trimmed -eq $oldMetaText) {
        $r = $para.Range
        $target = $d.Range($r.Start, $r.End - 1)
        $target.Text = $imagePromptText
    }
}
